# Auto-generated edit script applying recomputed market-price values
# to the Titan_Profits workbook (scheduled runner update).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 483595.97
$ws.Range("I28").Value = 617516.4399999999
$ws.Range("J28").Value = 1482.2
$ws.Range("K28").Value = 617516.4399999999
$ws.Range("L28").Value = 1482.2
$ws.Range("M28").Value = -617031.4399999999
$ws.Range("N28").Value = -2452.2
$ws.Range("H38").Value = 2062.7693
$ws.Range("J38").Value = 1828.375
$ws.Range("L38").Value = 5485.125
$ws.Range("N38").Value = -6229.125
$ws.Range("H43").Value = 583.3333
$ws.Range("J43").Value = 600
$ws.Range("L43").Value = 600
$ws.Range("N43").Value = -738
$ws.Range("H62").Value = 13374.25
$ws.Range("I62").Value = 8691.154
$ws.Range("J62").Value = 22071.428
$ws.Range("K62").Value = 8691.154
$ws.Range("L62").Value = 22071.428
$ws.Range("M62").Value = -8067.154
$ws.Range("N62").Value = -23319.428
$ws.Range("H65").Value = 13374.25
$ws.Range("I65").Value = 8691.154
$ws.Range("J65").Value = 22071.428
$ws.Range("K65").Value = 43455.77
$ws.Range("L65").Value = 110357.14
$ws.Range("M65").Value = -40335.77
$ws.Range("N65").Value = -116597.14
$ws.Range("H111").Value = 1422.125
$ws.Range("I111").Value = 1466.6666
$ws.Range("J111").Value = 1395.4
$ws.Range("K111").Value = 4399.9998
$ws.Range("L111").Value = 4186.200000000001
$ws.Range("M111").Value = -1332.9998
$ws.Range("N111").Value = -10320.2
$ws.Range("H116").Value = 4774075.5
$ws.Range("I116").Value = 6292172.5
$ws.Range("K116").Value = 6292172.5
$ws.Range("M116").Value = -6288730.5
$ws.Range("H132").Value = 217355.19
$ws.Range("I132").Value = 238832.36
$ws.Range("K132").Value = 716497.08
$ws.Range("M132").Value = -713967.08
$ws.Range("H137").Value = 35715468
$ws.Range("I137").Value = 40001016
$ws.Range("J137").Value = 2567.6667
$ws.Range("K137").Value = 120003048
$ws.Range("L137").Value = 7703.000100000001
$ws.Range("M137").Value = -120000498
$ws.Range("N137").Value = -12803.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 127527.25
$ws.Range("I2").Value = 202695.8
$ws.Range("J2").Value = 2246.3333
$ws.Range("K2").Value = 202695.8
$ws.Range("L2").Value = 2246.3333
$ws.Range("M2").Value = -202582.8
$ws.Range("N2").Value = -2472.3333
$ws.Range("H61").Value = 2020.8937
$ws.Range("I61").Value = 1300
$ws.Range("J61").Value = 4380.1816
$ws.Range("K61").Value = 1300
$ws.Range("L61").Value = 4380.1816
$ws.Range("M61").Value = -1088
$ws.Range("N61").Value = -4804.1816
$ws.Range("H110").Value = 1956.5
$ws.Range("I110").Value = 1300
$ws.Range("K110").Value = 1300
$ws.Range("M110").Value = 745
$ws.Range("H116").Value = 127527.25
$ws.Range("I116").Value = 202695.8
$ws.Range("J116").Value = 2246.3333
$ws.Range("K116").Value = 202695.8
$ws.Range("L116").Value = 2246.3333
$ws.Range("M116").Value = -200401.8
$ws.Range("N116").Value = -6834.3333
$ws.Range("H136").Value = 2020.8937
$ws.Range("I136").Value = 1300
$ws.Range("J136").Value = 4380.1816
$ws.Range("K136").Value = 3900
$ws.Range("L136").Value = 13140.5448
$ws.Range("M136").Value = -1350
$ws.Range("N136").Value = -18240.5448

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 127527.25
$ws.Range("I3").Value = 202695.8
$ws.Range("J3").Value = 2246.3333
$ws.Range("K3").Value = 202695.8
$ws.Range("L3").Value = 2246.3333
$ws.Range("M3").Value = -202581.8
$ws.Range("N3").Value = -2474.3333
$ws.Range("H92").Value = 4401
$ws.Range("J92").Value = 4401
$ws.Range("L92").Value = 4401
$ws.Range("N92").Value = -9393
$ws.Range("H105").Value = 3296.6667
$ws.Range("I105").Value = 3155.4443
$ws.Range("J105").Value = 3579.111
$ws.Range("K105").Value = 3155.4443
$ws.Range("L105").Value = 3579.111
$ws.Range("M105").Value = -1408.4443
$ws.Range("N105").Value = -7073.111
$ws.Range("H134").Value = 15875193
$ws.Range("I134").Value = 25001478
$ws.Range("J134").Value = 3392.739
$ws.Range("K134").Value = 75004434
$ws.Range("L134").Value = 10178.217
$ws.Range("M134").Value = -75001899
$ws.Range("N134").Value = -15248.217
$ws.Range("H135").Value = 75450.91
$ws.Range("J135").Value = 75450.91
$ws.Range("L135").Value = 75450.91
$ws.Range("N135").Value = -85590.91

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1448.3658
$ws.Range("I58").Value = 435.3793
$ws.Range("J58").Value = 3896.4167
$ws.Range("K58").Value = 435.3793
$ws.Range("L58").Value = 3896.4167
$ws.Range("M58").Value = -232.3793
$ws.Range("N58").Value = -4302.4167
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null
$ws.Range("H109").Value = 29733
$ws.Range("J109").Value = 29849.5
$ws.Range("L109").Value = 29849.5
$ws.Range("N109").Value = -31929.5
$ws.Range("H122").Value = 1734.375
$ws.Range("I122").Value = 1074.0667
$ws.Range("J122").Value = 2834.889
$ws.Range("K122").Value = 3222.2001
$ws.Range("L122").Value = 8504.667000000001
$ws.Range("M122").Value = -772.2001
$ws.Range("N122").Value = -13404.667
$ws.Range("H132").Value = 1983.0952
$ws.Range("I132").Value = 1701.1082
$ws.Range("J132").Value = 4069.8
$ws.Range("K132").Value = 5103.3246
$ws.Range("L132").Value = 12209.4
$ws.Range("M132").Value = -2573.3246
$ws.Range("N132").Value = -17269.4
$ws.Range("H136").Value = 1448.3658
$ws.Range("I136").Value = 435.3793
$ws.Range("J136").Value = 3896.4167
$ws.Range("K136").Value = 1306.1379
$ws.Range("L136").Value = 11689.2501
$ws.Range("M136").Value = 1243.8621
$ws.Range("N136").Value = -16789.2501

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 603.5714
$ws.Range("I122").Value = 269.75
$ws.Range("J122").Value = 1048.6666
$ws.Range("K122").Value = 2427.75
$ws.Range("L122").Value = 9437.999400000001
$ws.Range("M122").Value = 22.25
$ws.Range("N122").Value = -14337.9994
$ws.Range("H126").Value = 71431384
$ws.Range("I126").Value = 1250
$ws.Range("J126").Value = 83336410
$ws.Range("K126").Value = 3750
$ws.Range("L126").Value = 250009230
$ws.Range("M126").Value = 1190
$ws.Range("N126").Value = -250019110
$ws.Range("H130").Value = 995.8333
$ws.Range("J130").Value = 995.8333
$ws.Range("L130").Value = 2987.4999
$ws.Range("N130").Value = -13027.4999
$ws.Range("H131").Value = 1901.8064
$ws.Range("I131").Value = 432.125
$ws.Range("J131").Value = 2413
$ws.Range("K131").Value = 1296.375
$ws.Range("L131").Value = 7239
$ws.Range("M131").Value = 3743.625
$ws.Range("N131").Value = -17319

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 872.25
$ws.Range("I97").Value = 720.4167
$ws.Range("K97").Value = 720.4167
$ws.Range("M97").Value = -224.4167
$ws.Range("H113").Value = 1510
$ws.Range("I113").Value = 1013.3333
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1013.3333
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1156.6667
$ws.Range("N113").Value = -7340
$ws.Range("H126").Value = 3035.9412
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 3107.4
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 9322.200000000001
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -14262.2
$ws.Range("H132").Value = 4528.5654
$ws.Range("I132").Value = 4670.533
$ws.Range("J132").Value = 4262.375
$ws.Range("K132").Value = 14011.599
$ws.Range("L132").Value = 12787.125
$ws.Range("M132").Value = -11481.599
$ws.Range("N132").Value = -17847.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3435.2942
$ws.Range("I40").Value = 3200
$ws.Range("J40").Value = 3466.6667
$ws.Range("K40").Value = 3200
$ws.Range("L40").Value = 3466.6667
$ws.Range("M40").Value = -3064
$ws.Range("N40").Value = -3738.6667
$ws.Range("H46").Value = 2150
$ws.Range("I46").Value = 1150
$ws.Range("J46").Value = 3150
$ws.Range("K46").Value = 1150
$ws.Range("L46").Value = 3150
$ws.Range("M46").Value = -962
$ws.Range("N46").Value = -3526
$ws.Range("H48").Value = 50000
$ws.Range("I48").Value = 50000
$ws.Range("K48").Value = 50000
$ws.Range("M48").Value = -49339
$ws.Range("H55").Value = 465.41666
$ws.Range("J55").Value = 570.7143
$ws.Range("L55").Value = 570.7143
$ws.Range("N55").Value = -916.7143
$ws.Range("H100").Value = 3127539.2
$ws.Range("I100").Value = 12501617
$ws.Range("J100").Value = 2846.6667
$ws.Range("K100").Value = 12501617
$ws.Range("L100").Value = 2846.6667
$ws.Range("M100").Value = -12501076
$ws.Range("N100").Value = -3928.6667
$ws.Range("H136").Value = 4665.1943
$ws.Range("I136").Value = 2922.4443
$ws.Range("J136").Value = 9893.444
$ws.Range("K136").Value = 8767.332900000001
$ws.Range("L136").Value = 29680.332
$ws.Range("M136").Value = -6217.332900000001
$ws.Range("N136").Value = -34780.33199999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 256.5
$ws.Range("I113").Value = 254.57143
$ws.Range("J113").Value = 263.25
$ws.Range("K113").Value = 763.71429
$ws.Range("L113").Value = 789.75
$ws.Range("M113").Value = 1406.28571
$ws.Range("N113").Value = -5129.75
$ws.Range("H136").Value = 9037476
$ws.Range("I136").Value = 10785847
$ws.Range("K136").Value = 32357541
$ws.Range("M136").Value = -32354991

